$d = $word.ActiveDocument

# Fix typo: speigato -> spiegato
$d.Content.Find.Execute("speigato", $true, $false, $false, $false, $false, $true, 1, $false, "spiegato", 2)
